# technical_architecture.pptx - reposition several callout TextBoxes on slide 2.
# All affected shapes live inside the single top-level group shape ("그룹 106")
# that wraps the whole slide content, so we reach them via GroupItems.
#
# EMU -> point conversion uses the standard 12700 EMU-per-point factor. The
# host stores Shape.Left/Top/Width/Height as single-precision floats and
# truncates when converting back to EMU on save, so a plain x/12700.0 can
# land 1 EMU low. A tiny epsilon (well under half an EMU-in-points) nudges
# the float up enough to round-trip exactly without risking an off-by-one
# in the other direction.
$EMU_PER_PT = 12700.0
$EPS = 0.00004

function EmuToPt([double]$emu) {
    return ($emu / $EMU_PER_PT) + $EPS
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$grp = $s.Shapes.Item(1)

# Map of shape Id -> new Left/Top (and, where present, Width/Height) in EMU.
# $null entries mean "leave unchanged".
$updates = @{
    40  = @{ x = 4225051; y = 5664696 }                                   # TextBox 39 - 3.Varcode
    77  = @{ x = 7922465; y = 4646547 }                                   # TextBox 76 - 2.Varcode
    79  = @{ x = 9736560; y = 4666747 }                                   # TextBox 78 - 1.Varcode
    42  = @{ x = 4026680; y = 3886562; cx = 1206000; cy = 553998; wrap = $true }  # TextBox 41 - 5.Search/By Barcode
    44  = @{ x = 4115373; y = 4423883 }                                   # TextBox 43 - 6.Return/Item info
    48  = @{ x = 9935972; y = 7130855 }                                   # TextBox 47 - 1.x, y
    85  = @{ x = 7862257; y = 6240760 }                                   # TextBox 84 - 2.ax, ay, gz
    86  = @{ x = 8136820; y = 7069723 }                                   # TextBox 85 - 2.x, y
    49  = @{ x = 3247155; y = 6732808 }                                   # TextBox 48 - 3. Current Position & Direction
    94  = @{ x = 4112167; y = 2010877 }                                   # TextBox 93 - 2.Search/By Name
    100 = @{ x = 4115373; y = 2734434 }                                   # TextBox 99 - 3.Return/Item info
}

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $sh = $grp.GroupItems.Item($i)
    $u = $updates[$sh.Id]
    if ($u -ne $null) {
        if ($u.ContainsKey("cx")) {
            $sh.Width = EmuToPt $u.cx
        }
        if ($u.ContainsKey("cy")) {
            $sh.Height = EmuToPt $u.cy
        }
        $sh.Left = EmuToPt $u.x
        $sh.Top = EmuToPt $u.y
        if ($u.ContainsKey("wrap")) {
            $sh.TextFrame.WordWrap = -1
        }
    }
}
